$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Existing last row's timestamp got rewritten with slightly different
# floating-point precision by the logging tool on this run
$ws.Range("A70").Value2 = 44383.767590353

# New data row appended to the log (row 71)
$ws.Range("A71").NumberFormat = $ws.Range("A70").NumberFormat
$ws.Range("A71").Value2 = 44384.76911854006

$ws.Range("B71").Value = 79120
$ws.Range("C71").Value = 66782
$ws.Range("D71").Value = 3573
$ws.Range("E71").Value = 2162
$ws.Range("F71").Value = 1546
$ws.Range("G71").Value = 21075
$ws.Range("H71").Value = 1570
$ws.Range("I71").Value = 885
$ws.Range("J71").Value = 203
